# Automatische test-sync: 2025-07-31 21:43:50
# Appends testmail #10 ("Is er al nieuws?") as a new row to the Logs sheet,
# extends the conditional formatting ranges to cover the new row, and
# bumps the "Overig" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 12

$ws.Cells.Item($newRow, 1).Value = "Is er al nieuws?"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #10: Is er al nieuws?"
$ws.Cells.Item($newRow, 4).Value = "Overig"
$ws.Cells.Item($newRow, 5).Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Cells.Item($newRow, 6).Value = "2025-07-31 21:43:42"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting ranges (D, G, H, I, J) from row 11 to row 12.
# All cfRules that share the same sqref move together when one of them is
# re-applied, so touching the first rule in each column group is enough.
$ws.Range("D2:D11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D12"))
$ws.Range("G2:G11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G12"))
$ws.Range("H2:H11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H12"))
$ws.Range("I2:I11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I12"))
$ws.Range("J2:J11").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J12"))

# Update the Dashboard summary: "Overig" count goes from 4 to 5.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 5
